# Update cryptos list: price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'58.751.82"
$ws.Range('E2').Value = '  -2.34%  '
$ws.Range('D3').Value = "'2.293.50"
$ws.Range('E3').Value = '  -5.36%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'546.74"
$ws.Range('E5').Value = '  -1.24%  '
$ws.Range('D6').Value = "'130.60"
$ws.Range('E6').Value = '  -4.85%  '
$ws.Range('D8').Value = "'0.570"
$ws.Range('E8').Value = '  -3.20%  '
$ws.Range('D9').Value = "'2.292.02"
$ws.Range('E9').Value = '  -5.37%  '
$ws.Range('E10').Value = '  -3.50%  '
$ws.Range('E11').Value = '  -2.99%  '
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('E13').Value = '  -5.23%  '
$ws.Range('D14').Value = "'23.83"
$ws.Range('E14').Value = '  -4.31%  '
$ws.Range('D15').Value = "'2.699.70"
$ws.Range('E15').Value = '  -5.46%  '
$ws.Range('D16').Value = "'58.698.81"
$ws.Range('E16').Value = '  -2.31%  '
$ws.Range('E17').Value = '  -3.59%  '
$ws.Range('D18').Value = "'2.348.37"
$ws.Range('E18').Value = '  -3.05%  '
$ws.Range('D19').Value = "'10.66"
$ws.Range('E19').Value = '  -5.51%  '
$ws.Range('E20').Value = '  -4.42%  '
$ws.Range('D21').Value = "'315.05"
$ws.Range('E21').Value = '  -3.89%  '
$ws.Range('E22').Value = '  -4.15%  '
$ws.Range('D23').Value = "'1.00"
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').Value = "'63.08"
$ws.Range('E24').Value = '  -3.23%  '
$ws.Range('D25').Value = "'0.172"
$ws.Range('E25').Value = '  -3.24%  '
$ws.Range('D26').Value = "'0.997"
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('E27').Value = '  -6.85%  '
$ws.Range('D28').Value = "'1.31"
$ws.Range('E28').Value = '  -6.60%  '
$ws.Range('E29').Value = '  -1.98%  '
$ws.Range('D30').Value = "'169.45"
$ws.Range('E30').Value = '  -0.66%  '
$ws.Range('D31').Value = "'0.0₃0728"
$ws.Range('E31').Value = '  -5.97%  '
$ws.Range('D32').Value = "'5.79"
$ws.Range('E32').Value = '  -5.28%  '
$ws.Range('D33').Value = "'1.08"
$ws.Range('E33').Value = '  +1.20%  '
$ws.Range('D34').Value = "'0.385"
$ws.Range('E34').Value = '  -4.89%  '
$ws.Range('D36').Value = "'17.81"
$ws.Range('E36').Value = '  -4.12%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = "'1.25"
$ws.Range('E38').Value = '  -6.75%  '
$ws.Range('D39').Value = "'3.96"
$ws.Range('E39').Value = '  -6.45%  '
$ws.Range('D40').Value = "'37.95"
$ws.Range('E41').Value = '  -5.65%  '
$ws.Range('D42').Value = "'301.70"
$ws.Range('E42').Value = '  -7.33%  '
$ws.Range('D43').Value = "'140.02"
$ws.Range('E43').Value = '  -3.49%  '
$ws.Range('D44').Value = "'3.45"
$ws.Range('E44').Value = '  -5.78%  '
$ws.Range('D45').Value = "'0.0951"
$ws.Range('E45').Value = '  -1.37%  '
$ws.Range('D46').Value = "'0.0499"
$ws.Range('E46').Value = '  -3.69%  '
$ws.Range('E47').Value = '  -3.28%  '
$ws.Range('E48').Value = '  -7.23%  '
$ws.Range('E49').Value = '  -3.76%  '
$ws.Range('D50').Value = "'16.63"
$ws.Range('E50').Value = '  -5.07%  '
$ws.Range('E51').Value = '  -0.39%  '
